$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 5500.5
$ws.Range("I12").Value = 5500.5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 5500.5
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = $null
$ws.Range("N12").Value = -5330.5
$ws.Range("H29").Value = 1582.7084
$ws.Range("J29").Value = 1513.5714
$ws.Range("L29").Value = 4540.7142
$ws.Range("N29").Value = -5102.7142
$ws.Range("H70").Value = 4168.857
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 4530.3335
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 13591.0005
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = -14131.0005
$ws.Range("H73").Value = 4168.857
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 4530.3335
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 13591.0005
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = -15463.0005
$ws.Range("H87").Value = 58500
$ws.Range("J87").Value = 58500
$ws.Range("L87").Value = 58500
$ws.Range("N87").Value = -60996
$ws.Range("H90").Value = 58500
$ws.Range("J90").Value = 58500
$ws.Range("L90").Value = 175500
$ws.Range("N90").Value = -187980
$ws.Range("H92").Value = 1434.9375
$ws.Range("I92").Value = 2583.125
$ws.Range("J92").Value = 286.75
$ws.Range("K92").Value = 2583.125
$ws.Range("L92").Value = 286.75
$ws.Range("M92").Value = -1335.125
$ws.Range("N92").Value = -2782.75
$ws.Range("H96").Value = 1463.4166
$ws.Range("I96").Value = 337.25
$ws.Range("J96").Value = 2589.5833
$ws.Range("K96").Value = 1011.75
$ws.Range("L96").Value = 7768.749899999999
$ws.Range("M96").Value = 361.25
$ws.Range("N96").Value = -10514.7499
$ws.Range("H112").Value = 1777.4117
$ws.Range("J112").Value = 1786.4615
$ws.Range("L112").Value = 5359.3845
$ws.Range("N112").Value = -7575.3845
$ws.Range("H132").Value = 14524.059
$ws.Range("I132").Value = 14524.059
$ws.Range("K132").Value = 43572.177
$ws.Range("M132").Value = -41042.177
$ws.Range("H137").Value = 3576907.2
$ws.Range("I137").Value = 6252010
$ws.Range("J137").Value = 10103.625
$ws.Range("K137").Value = 18756030
$ws.Range("L137").Value = 30310.875
$ws.Range("M137").Value = -18753480
$ws.Range("N137").Value = -35410.875
$ws.Range("H138").Value = 6311.459
$ws.Range("I138").Value = 6516.7646
$ws.Range("J138").Value = 6232.136
$ws.Range("K138").Value = 19550.2938
$ws.Range("L138").Value = 18696.408
$ws.Range("M138").Value = -14410.2938
$ws.Range("N138").Value = -28976.408
$ws.Range("H141").Value = 4817.5386
$ws.Range("I141").Value = 1763.15
$ws.Range("K141").Value = 5289.450000000001
$ws.Range("M141").Value = -109.4500000000007

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2769.6667
$ws.Range("I45").Value = 2463.7334
$ws.Range("J45").Value = 4299.3335
$ws.Range("K45").Value = 2463.7334
$ws.Range("L45").Value = 4299.3335
$ws.Range("M45").Value = -2086.7334
$ws.Range("N45").Value = -5053.3335
$ws.Range("H61").Value = 5927912
$ws.Range("I61").Value = 3970217.5
$ws.Range("J61").Value = 33335634
$ws.Range("K61").Value = 3970217.5
$ws.Range("L61").Value = 33335634
$ws.Range("M61").Value = -3970005.5
$ws.Range("N61").Value = -33336058
$ws.Range("H97").Value = 1556.4615
$ws.Range("I97").Value = 317.44446
$ws.Range("J97").Value = 4344.25
$ws.Range("K97").Value = 317.44446
$ws.Range("L97").Value = 4344.25
$ws.Range("M97").Value = 178.55554
$ws.Range("N97").Value = -5336.25
$ws.Range("H102").Value = 2603.4546
$ws.Range("I102").Value = 2713.8
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 2713.8
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = -1091.8
$ws.Range("N102").Value = -4744
$ws.Range("H122").Value = 35716228
$ws.Range("I122").Value = 41668430
$ws.Range("J122").Value = 2986.5
$ws.Range("K122").Value = 125005290
$ws.Range("L122").Value = 8959.5
$ws.Range("M122").Value = -125002840
$ws.Range("N122").Value = -13859.5
$ws.Range("H136").Value = 5927912
$ws.Range("I136").Value = 3970217.5
$ws.Range("J136").Value = 33335634
$ws.Range("K136").Value = 11910652.5
$ws.Range("L136").Value = 100006902
$ws.Range("M136").Value = -11908102.5
$ws.Range("N136").Value = -100012002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1738.7778
$ws.Range("I86").Value = 1388.7693
$ws.Range("J86").Value = 2648.8
$ws.Range("K86").Value = 1388.7693
$ws.Range("L86").Value = 2648.8
$ws.Range("M86").Value = -265.7692999999999
$ws.Range("N86").Value = -4894.8
$ws.Range("H89").Value = 1738.7778
$ws.Range("I89").Value = 1388.7693
$ws.Range("J89").Value = 2648.8
$ws.Range("K89").Value = 6943.8465
$ws.Range("L89").Value = 13244
$ws.Range("M89").Value = -1327.8465
$ws.Range("N89").Value = -24476
$ws.Range("H94").Value = 2773
$ws.Range("I94").Value = 3244.762
$ws.Range("J94").Value = 296.25
$ws.Range("K94").Value = 3244.762
$ws.Range("L94").Value = 296.25
$ws.Range("M94").Value = -2793.762
$ws.Range("N94").Value = -1198.25
$ws.Range("H99").Value = 13080.286
$ws.Range("I99").Value = 14374.909
$ws.Range("J99").Value = 8333.333000000001
$ws.Range("K99").Value = 14374.909
$ws.Range("L99").Value = 8333.333000000001
$ws.Range("M99").Value = -12876.909
$ws.Range("N99").Value = -11329.333
$ws.Range("H134").Value = 12094779
$ws.Range("I134").Value = 14515325
$ws.Range("J134").Value = 6951117.5
$ws.Range("K134").Value = 43545975
$ws.Range("L134").Value = 20853352.5
$ws.Range("M134").Value = -43543440
$ws.Range("N134").Value = -20858422.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1719.9565
$ws.Range("I22").Value = 1170.8667
$ws.Range("J22").Value = 2749.5
$ws.Range("K22").Value = 1170.8667
$ws.Range("L22").Value = 2749.5
$ws.Range("M22").Value = -820.8667
$ws.Range("N22").Value = -3449.5
$ws.Range("H31").Value = 435272.1
$ws.Range("I31").Value = 801249.0600000001
$ws.Range("J31").Value = 3942.1072
$ws.Range("K31").Value = 801249.0600000001
$ws.Range("L31").Value = 3942.1072
$ws.Range("M31").Value = -800954.0600000001
$ws.Range("N31").Value = -4532.1072
$ws.Range("H34").Value = 435272.1
$ws.Range("I34").Value = 801249.0600000001
$ws.Range("J34").Value = 3942.1072
$ws.Range("K34").Value = 801249.0600000001
$ws.Range("L34").Value = 3942.1072
$ws.Range("M34").Value = -801047.0600000001
$ws.Range("N34").Value = -4346.1072
$ws.Range("H134").Value = 4501.9346
$ws.Range("I134").Value = 2783.0557
$ws.Range("J134").Value = 5606.9287
$ws.Range("K134").Value = 8349.167099999999
$ws.Range("L134").Value = 16820.7861
$ws.Range("M134").Value = -5814.167099999999
$ws.Range("N134").Value = -21890.7861

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 16500968
$ws.Range("I4").Value = 21000242
$ws.Range("J4").Value = 3633
$ws.Range("K4").Value = 63000726
$ws.Range("L4").Value = 10899
$ws.Range("M4").Value = -63000614
$ws.Range("N4").Value = -11123
$ws.Range("H7").Value = 36.666668
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = $null
$ws.Range("H57").Value = 999997
$ws.Range("I57").Value = 999997
$ws.Range("K57").Value = 2999991
$ws.Range("M57").Value = -2999432
$ws.Range("H74").Value = 14166.667
$ws.Range("H76").Value = 4999
$ws.Range("J76").Value = 4999
$ws.Range("L76").Value = 14997
$ws.Range("N76").Value = -15763
$ws.Range("H77").Value = 14166.667
$ws.Range("H79").Value = 4999
$ws.Range("J79").Value = 4999
$ws.Range("L79").Value = 14997
$ws.Range("N79").Value = -17649
$ws.Range("H80").Value = 15129.223
$ws.Range("J80").Value = 16845.375
$ws.Range("L80").Value = 50536.125
$ws.Range("N80").Value = -52408.125
$ws.Range("H83").Value = 15129.223
$ws.Range("J83").Value = 16845.375
$ws.Range("L83").Value = 151608.375
$ws.Range("N83").Value = -160968.375
$ws.Range("H113").Value = 568
$ws.Range("I113").Value = 429.45456
$ws.Range("K113").Value = 1288.36368
$ws.Range("M113").Value = 881.6363200000001
$ws.Range("H122").Value = 646097.5
$ws.Range("I122").Value = 1152657.4
$ws.Range("J122").Value = 1384.8182
$ws.Range("K122").Value = 10373916.6
$ws.Range("L122").Value = 12463.3638
$ws.Range("M122").Value = -10371466.6
$ws.Range("N122").Value = -17363.3638
$ws.Range("H140").Value = 18532.695
$ws.Range("I140").Value = 19147.818
$ws.Range("K140").Value = 57443.454
$ws.Range("M140").Value = -52263.454

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 952.8570999999999
$ws.Range("I97").Value = 870.2414
$ws.Range("J97").Value = 1137.1538
$ws.Range("K97").Value = 870.2414
$ws.Range("L97").Value = 1137.1538
$ws.Range("M97").Value = -374.2414
$ws.Range("N97").Value = -2129.1538
$ws.Range("H102").Value = 3404.5
$ws.Range("I102").Value = 2202.1035
$ws.Range("J102").Value = 10378.4
$ws.Range("K102").Value = 2202.1035
$ws.Range("L102").Value = 10378.4
$ws.Range("M102").Value = -580.1035000000002
$ws.Range("N102").Value = -13622.4
$ws.Range("H113").Value = 2541.4666
$ws.Range("I113").Value = 2194
$ws.Range("K113").Value = 2194
$ws.Range("M113").Value = -24
$ws.Range("H122").Value = 6784969.5
$ws.Range("I122").Value = 8381165
$ws.Range("J122").Value = 1137
$ws.Range("K122").Value = 25143495
$ws.Range("L122").Value = 3411
$ws.Range("M122").Value = -25141045
$ws.Range("N122").Value = -8311

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 11042.857
$ws.Range("I122").Value = 12560.4
$ws.Range("J122").Value = 7249
$ws.Range("K122").Value = 37681.2
$ws.Range("L122").Value = 21747
$ws.Range("M122").Value = -35231.2
$ws.Range("N122").Value = -26647

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1749.5
$ws.Range("I96").Value = 1749.5
$ws.Range("K96").Value = 1749.5
$ws.Range("M96").Value = -376.5
$ws.Range("H113").Value = 1560.2188
$ws.Range("I113").Value = 1060.2
$ws.Range("J113").Value = 2001.4117
$ws.Range("K113").Value = 3180.6
$ws.Range("L113").Value = 6004.2351
$ws.Range("M113").Value = -1010.6
$ws.Range("N113").Value = -10344.2351
$ws.Range("H122").Value = 71211.17999999999
$ws.Range("I122").Value = 5584.6
$ws.Range("K122").Value = 16753.8
$ws.Range("M122").Value = -14303.8
$ws.Range("H126").Value = 2424.2144
$ws.Range("I126").Value = 2649.3333
$ws.Range("J126").Value = 2019
$ws.Range("K126").Value = 7947.999899999999
$ws.Range("L126").Value = 6057
$ws.Range("M126").Value = -5477.999899999999
$ws.Range("N126").Value = -10997
